# Fix malformed/quoted strings in the "purpose" column (C) and restore
# the quoted "job" value that lost its quotes (H418, H733).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "new car"
$ws.Range("C9").Value = "used car"
$ws.Range("C11").Value = "new car"
$ws.Range("C12").Value = "new car"
$ws.Range("C15").Value = "new car"
$ws.Range("C16").Value = "new car"
$ws.Range("C20").Value = "used car"
$ws.Range("C22").Value = "new car"
$ws.Range("C24").Value = "new car"
$ws.Range("C25").Value = "used car"
$ws.Range("C34").Value = "new car"
$ws.Range("C40").Value = "domestic appliance"
$ws.Range("C45").Value = "used car"
$ws.Range("C46").Value = "used car"
$ws.Range("C47").Value = "new car"
$ws.Range("C49").Value = "used car"
$ws.Range("C50").Value = "new car"
$ws.Range("C53").Value = "used car"
$ws.Range("C55").Value = "used car"
$ws.Range("C56").Value = "new car"
$ws.Range("C57").Value = "new car"
$ws.Range("C60").Value = "new car"
$ws.Range("C69").Value = "new car"
$ws.Range("C72").Value = "used car"
$ws.Range("C77").Value = "used car"
$ws.Range("C80").Value = "used car"
$ws.Range("C85").Value = "other"
$ws.Range("C90").Value = "new car"
$ws.Range("C93").Value = "used car"
$ws.Range("C96").Value = "new car"
$ws.Range("C101").Value = "used car"
$ws.Range("C102").Value = "new car"
$ws.Range("C106").Value = "used car"
$ws.Range("C108").Value = "new car"
$ws.Range("C109").Value = "new car"
$ws.Range("C114").Value = "new car"
$ws.Range("C115").Value = "new car"
$ws.Range("C121").Value = "new car"
$ws.Range("C123").Value = "used car"
$ws.Range("C125").Value = "new car"
$ws.Range("C127").Value = "new car"
$ws.Range("C130").Value = "used car"
$ws.Range("C131").Value = "new car"
$ws.Range("C132").Value = "new car"
$ws.Range("C138").Value = "used car"
$ws.Range("C142").Value = "new car"
$ws.Range("C148").Value = "new car"
$ws.Range("C149").Value = "new car"
$ws.Range("C155").Value = "used car"
$ws.Range("C160").Value = "new car"
$ws.Range("C163").Value = "new car"
$ws.Range("C164").Value = "domestic appliance"
$ws.Range("C165").Value = "new car"
$ws.Range("C166").Value = "new car"
$ws.Range("C172").Value = "new car"
$ws.Range("C177").Value = "used car"
$ws.Range("C181").Value = "new car"
$ws.Range("C184").Value = "new car"
$ws.Range("C186").Value = "new car"
$ws.Range("C188").Value = "used car"
$ws.Range("C189").Value = "new car"
$ws.Range("C200").Value = "used car"
$ws.Range("C203").Value = "new car"
$ws.Range("C206").Value = "new car"
$ws.Range("C207").Value = "used car"
$ws.Range("C209").Value = "domestic appliance"
$ws.Range("C211").Value = "used car"
$ws.Range("C221").Value = "new car"
$ws.Range("C226").Value = "used car"
$ws.Range("C233").Value = "new car"
$ws.Range("C238").Value = "new car"
$ws.Range("C242").Value = "new car"
$ws.Range("C244").Value = "used car"
$ws.Range("C249").Value = "new car"
$ws.Range("C252").Value = "new car"
$ws.Range("C254").Value = "new car"
$ws.Range("C260").Value = "used car"
$ws.Range("C264").Value = "new car"
$ws.Range("C266").Value = "new car"
$ws.Range("C270").Value = "new car"
$ws.Range("C272").Value = "new car"
$ws.Range("C274").Value = "new car"
$ws.Range("C282").Value = "used car"
$ws.Range("C286").Value = "new car"
$ws.Range("C287").Value = "new car"
$ws.Range("C288").Value = "used car"
$ws.Range("C293").Value = "used car"
$ws.Range("C294").Value = "used car"
$ws.Range("C295").Value = "used car"
$ws.Range("C298").Value = "used car"
$ws.Range("C299").Value = "new car"
$ws.Range("C302").Value = "new car"
$ws.Range("C304").Value = "new car"
$ws.Range("C305").Value = "new car"
$ws.Range("C306").Value = "new car"
$ws.Range("C308").Value = "used car"
$ws.Range("C311").Value = "new car"
$ws.Range("C315").Value = "new car"
$ws.Range("C316").Value = "new car"
$ws.Range("C322").Value = "new car"
$ws.Range("C324").Value = "used car"
$ws.Range("C326").Value = "new car"
$ws.Range("C327").Value = "new car"
$ws.Range("C329").Value = "new car"
$ws.Range("C332").Value = "used car"
$ws.Range("C334").Value = "new car"
$ws.Range("C335").Value = "used car"
$ws.Range("C339").Value = "domestic appliance"
$ws.Range("C346").Value = "new car"
$ws.Range("C352").Value = "domestic appliance"
$ws.Range("C354").Value = "used car"
$ws.Range("C357").Value = "new car"
$ws.Range("C364").Value = "new car"
$ws.Range("C368").Value = "used car"
$ws.Range("C372").Value = "new car"
$ws.Range("C375").Value = "new car"
$ws.Range("C376").Value = "other"
$ws.Range("C380").Value = "new car"
$ws.Range("C381").Value = "new car"
$ws.Range("C383").Value = "used car"
$ws.Range("C384").Value = "new car"
$ws.Range("C385").Value = "new car"
$ws.Range("C392").Value = "new car"
$ws.Range("C394").Value = "new car"
$ws.Range("C400").Value = "new car"
$ws.Range("C405").Value = "new car"
$ws.Range("C406").Value = "new car"
$ws.Range("C408").Value = "used car"
$ws.Range("C411").Value = "new car"
$ws.Range("C413").Value = "used car"
$ws.Range("C415").Value = "new car"
$ws.Range("C416").Value = "new car"
$ws.Range("C417").Value = "used car"
$ws.Range("C418").Value = "new car"
$ws.Range("C420").Value = "new car"
$ws.Range("C421").Value = "new car"
$ws.Range("C422").Value = "new car"
$ws.Range("C423").Value = "used car"
$ws.Range("C424").Value = "new car"
$ws.Range("C427").Value = "used car"
$ws.Range("C439").Value = "new car"
$ws.Range("C442").Value = "new car"
$ws.Range("C448").Value = "new car"
$ws.Range("C452").Value = "used car"
$ws.Range("C455").Value = "used car"
$ws.Range("C456").Value = "new car"
$ws.Range("C457").Value = "used car"
$ws.Range("C458").Value = "new car"
$ws.Range("C459").Value = "used car"
$ws.Range("C460").Value = "domestic appliance"
$ws.Range("C463").Value = "new car"
$ws.Range("C467").Value = "used car"
$ws.Range("C471").Value = "used car"
$ws.Range("C474").Value = "new car"
$ws.Range("C478").Value = "used car"
$ws.Range("C483").Value = "new car"
$ws.Range("C486").Value = "new car"
$ws.Range("C487").Value = "new car"
$ws.Range("C489").Value = "new car"
$ws.Range("C490").Value = "new car"
$ws.Range("C491").Value = "new car"
$ws.Range("C496").Value = "new car"
$ws.Range("C501").Value = "new car"
$ws.Range("C502").Value = "new car"
$ws.Range("C503").Value = "used car"
$ws.Range("C506").Value = "new car"
$ws.Range("C507").Value = "new car"
$ws.Range("C508").Value = "used car"
$ws.Range("C509").Value = "new car"
$ws.Range("C511").Value = "used car"
$ws.Range("C512").Value = "new car"
$ws.Range("C513").Value = "used car"
$ws.Range("C516").Value = "new car"
$ws.Range("C517").Value = "new car"
$ws.Range("C518").Value = "new car"
$ws.Range("C520").Value = "new car"
$ws.Range("C525").Value = "used car"
$ws.Range("C527").Value = "used car"
$ws.Range("C531").Value = "new car"
$ws.Range("C533").Value = "new car"
$ws.Range("C534").Value = "used car"
$ws.Range("C538").Value = "new car"
$ws.Range("C540").Value = "new car"
$ws.Range("C543").Value = "new car"
$ws.Range("C546").Value = "new car"
$ws.Range("C547").Value = "new car"
$ws.Range("C548").Value = "new car"
$ws.Range("C551").Value = "used car"
$ws.Range("C555").Value = "new car"
$ws.Range("C558").Value = "new car"
$ws.Range("C559").Value = "new car"
$ws.Range("C562").Value = "used car"
$ws.Range("C565").Value = "new car"
$ws.Range("C568").Value = "new car"
$ws.Range("C574").Value = "used car"
$ws.Range("C580").Value = "new car"
$ws.Range("C582").Value = "new car"
$ws.Range("C583").Value = "new car"
$ws.Range("C586").Value = "new car"
$ws.Range("C588").Value = "new car"
$ws.Range("C590").Value = "domestic appliance"
$ws.Range("C593").Value = "new car"
$ws.Range("C595").Value = "new car"
$ws.Range("C597").Value = "new car"
$ws.Range("C598").Value = "new car"
$ws.Range("C600").Value = "new car"
$ws.Range("C611").Value = "used car"
$ws.Range("C612").Value = "domestic appliance"
$ws.Range("C613").Value = "new car"
$ws.Range("C615").Value = "used car"
$ws.Range("C619").Value = "new car"
$ws.Range("C621").Value = "new car"
$ws.Range("C623").Value = "new car"
$ws.Range("C629").Value = "new car"
$ws.Range("C636").Value = "new car"
$ws.Range("C649").Value = "new car"
$ws.Range("C650").Value = "new car"
$ws.Range("C654").Value = "new car"
$ws.Range("C655").Value = "new car"
$ws.Range("C656").Value = "used car"
$ws.Range("C657").Value = "new car"
$ws.Range("C658").Value = "new car"
$ws.Range("C663").Value = "new car"
$ws.Range("C670").Value = "new car"
$ws.Range("C674").Value = "new car"
$ws.Range("C675").Value = "new car"
$ws.Range("C687").Value = "new car"
$ws.Range("C689").Value = "new car"
$ws.Range("C691").Value = "new car"
$ws.Range("C695").Value = "new car"
$ws.Range("C697").Value = "used car"
$ws.Range("C699").Value = "new car"
$ws.Range("C703").Value = "used car"
$ws.Range("C707").Value = "new car"
$ws.Range("C708").Value = "new car"
$ws.Range("C714").Value = "used car"
$ws.Range("C716").Value = "new car"
$ws.Range("C717").Value = "used car"
$ws.Range("C721").Value = "used car"
$ws.Range("C724").Value = "new car"
$ws.Range("C726").Value = "new car"
$ws.Range("C727").Value = "new car"
$ws.Range("C735").Value = "used car"
$ws.Range("C736").Value = "new car"
$ws.Range("C737").Value = "domestic appliance"
$ws.Range("C738").Value = "used car"
$ws.Range("C739").Value = "new car"
$ws.Range("C740").Value = "new car"
$ws.Range("C742").Value = "new car"
$ws.Range("C748").Value = "new car"
$ws.Range("C749").Value = "new car"
$ws.Range("C750").Value = "used car"
$ws.Range("C751").Value = "used car"
$ws.Range("C753").Value = "new car"
$ws.Range("C757").Value = "new car"
$ws.Range("C758").Value = "new car"
$ws.Range("C760").Value = "new car"
$ws.Range("C761").Value = "new car"
$ws.Range("C762").Value = "new car"
$ws.Range("C765").Value = "new car"
$ws.Range("C766").Value = "new car"
$ws.Range("C769").Value = "used car"
$ws.Range("C772").Value = "used car"
$ws.Range("C774").Value = "used car"
$ws.Range("C776").Value = "new car"
$ws.Range("C777").Value = "new car"
$ws.Range("C778").Value = "new car"
$ws.Range("C780").Value = "used car"
$ws.Range("C783").Value = "new car"
$ws.Range("C785").Value = "new car"
$ws.Range("C786").Value = "used car"
$ws.Range("C789").Value = "used car"
$ws.Range("C793").Value = "used car"
$ws.Range("C798").Value = "used car"
$ws.Range("C800").Value = "new car"
$ws.Range("C801").Value = "new car"
$ws.Range("C806").Value = "new car"
$ws.Range("C807").Value = "new car"
$ws.Range("C810").Value = "used car"
$ws.Range("C811").Value = "new car"
$ws.Range("C814").Value = "used car"
$ws.Range("C815").Value = "domestic appliance"
$ws.Range("C816").Value = "new car"
$ws.Range("C817").Value = "new car"
$ws.Range("C818").Value = "domestic appliance"
$ws.Range("C822").Value = "new car"
$ws.Range("C825").Value = "new car"
$ws.Range("C827").Value = "new car"
$ws.Range("C828").Value = "new car"
$ws.Range("C830").Value = "used car"
$ws.Range("C833").Value = "new car"
$ws.Range("C837").Value = "new car"
$ws.Range("C840").Value = "used car"
$ws.Range("C843").Value = "used car"
$ws.Range("C848").Value = "new car"
$ws.Range("C849").Value = "new car"
$ws.Range("C852").Value = "new car"
$ws.Range("C853").Value = "used car"
$ws.Range("C855").Value = "new car"
$ws.Range("C856").Value = "new car"
$ws.Range("C857").Value = "new car"
$ws.Range("C860").Value = "new car"
$ws.Range("C861").Value = "new car"
$ws.Range("C862").Value = "used car"
$ws.Range("C873").Value = "new car"
$ws.Range("C875").Value = "domestic appliance"
$ws.Range("C877").Value = "new car"
$ws.Range("C880").Value = "new car"
$ws.Range("C882").Value = "used car"
$ws.Range("C883").Value = "used car"
$ws.Range("C884").Value = "new car"
$ws.Range("C890").Value = "new car"
$ws.Range("C891").Value = "used car"
$ws.Range("C894").Value = "new car"
$ws.Range("C895").Value = "used car"
$ws.Range("C897").Value = "used car"
$ws.Range("C902").Value = "new car"
$ws.Range("C903").Value = "new car"
$ws.Range("C904").Value = "used car"
$ws.Range("C908").Value = "new car"
$ws.Range("C910").Value = "used car"
$ws.Range("C911").Value = "new car"
$ws.Range("C918").Value = "used car"
$ws.Range("C919").Value = "new car"
$ws.Range("C925").Value = "new car"
$ws.Range("C927").Value = "new car"
$ws.Range("C929").Value = "used car"
$ws.Range("C931").Value = "new car"
$ws.Range("C934").Value = "new car"
$ws.Range("C941").Value = "used car"
$ws.Range("C942").Value = "new car"
$ws.Range("C943").Value = "new car"
$ws.Range("C945").Value = "new car"
$ws.Range("C947").Value = "new car"
$ws.Range("C949").Value = "new car"
$ws.Range("C954").Value = "used car"
$ws.Range("C956").Value = "new car"
$ws.Range("C960").Value = "new car"
$ws.Range("C963").Value = "new car"
$ws.Range("C964").Value = "new car"
$ws.Range("C971").Value = "new car"
$ws.Range("C973").Value = "new car"
$ws.Range("C974").Value = "new car"
$ws.Range("C980").Value = "new car"
$ws.Range("C981").Value = "new car"
$ws.Range("C984").Value = "new car"
$ws.Range("C985").Value = "used car"
$ws.Range("C990").Value = "used car"
$ws.Range("C996").Value = "new car"
$ws.Range("C998").Value = "used car"
$ws.Range("C1001").Value = "used car"

# Restore the leading/trailing single quotes around these two cells.
# A literal leading quote is doubled so Excel's "force text" prefix
# strip only consumes the first one, leaving the real quote intact.
$ws.Range("H418").Value = "''unskilled resident'"
$ws.Range("H733").Value = "''unskilled resident'"

Write-Host "Applied purpose/job text fixes"
